$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 72: mdfilegenerator (merged O72:Q72) ---
$ws.Range("M72").Value = 44685
$ws.Range("M72").NumberFormat = "m/d/yyyy"

$ws.Range("O71:Q71").Copy()
$ws.Range("O72:Q72").PasteSpecial(-4122)
$ws.Range("O72:Q72").Merge()
$ws.Range("O72").Value = "mdfilegenerator"

$ws.Range("R72").Value = 0.66666666666666663
$ws.Range("R72").NumberFormat = "h:mm"
$ws.Range("S72").Value = 0.75
$ws.Range("S72").NumberFormat = "h:mm"
$ws.Range("T72").Formula = "=S72-R72"
$ws.Range("T72").NumberFormat = "h:mm"

# --- Row 73: Improving mdfilegenerator ---
$ws.Range("M73").Value = 44723
$ws.Range("M73").NumberFormat = "m/d/yyyy"

$ws.Range("O73").Value = "Improving mdfilegenerator"

$ws.Range("R73").Value = 0.70833333333333337
$ws.Range("R73").NumberFormat = "h:mm"
$ws.Range("S73").Value = 0.82291666666666663
$ws.Range("S73").NumberFormat = "h:mm"
$ws.Range("T73").Formula = "=S73-R73"
$ws.Range("T73").NumberFormat = "h:mm"

# --- Row 74: Working on mdfilegenerator ---
$ws.Range("M74").Value = 44724
$ws.Range("M74").NumberFormat = "m/d/yyyy"

$ws.Range("O74").Value = "Working on mdfilegenerator"

$ws.Range("R74").Value = 0.625
$ws.Range("R74").NumberFormat = "h:mm"
$ws.Range("S74").Value = 0.75
$ws.Range("S74").NumberFormat = "h:mm"
$ws.Range("T74").Formula = "=S74-R74"
$ws.Range("T74").NumberFormat = "h:mm"

# --- Row 75: Website frontend prüfen für finalen Sprint ---
$ws.Range("M75").Value = 44724
$ws.Range("M75").NumberFormat = "m/d/yyyy"

$ws.Range("O75").Value = "Website frontend prüfen für finalen Sprint"

$ws.Range("R75").Value = 0.85416666666666663
$ws.Range("R75").NumberFormat = "h:mm"
$ws.Range("S75").Value = 0.91666666666666663
$ws.Range("S75").NumberFormat = "h:mm"
$ws.Range("T75").Formula = "=S75-R75"
$ws.Range("T75").NumberFormat = "h:mm"

# --- Row 76: Powerpoint prüfen(Rechtschreibung, Grammatik, Content) ---
$ws.Range("M76").Value = 44724
$ws.Range("M76").NumberFormat = "m/d/yyyy"

$ws.Range("O76").Value = "Powerpoint prüfen(Rechtschreibung, Grammatik, Content)"

$ws.Range("R76").Value = 0.91666666666666663
$ws.Range("R76").NumberFormat = "h:mm"
$ws.Range("S76").Value = 0.95833333333333337
$ws.Range("S76").NumberFormat = "h:mm"
$ws.Range("T76").Formula = "=S76-R76"
$ws.Range("T76").NumberFormat = "h:mm"

# --- Fix summary formulas in H17 / H18 ---
$ws.Range("H17").Formula = "=SUM(T71:T76)"
$ws.Range("H18").Formula = "=SUM(AD71:AD76)"

# --- Extend sheet with a trailing empty row 78 (mirrors row 77) ---
$ws.Range("T77").Copy()
$ws.Range("T78").PasteSpecial(-4122)

Write-Output "done"
